$d = $word.ActiveDocument

# Locate the paragraph that ends the "Create Attendance Sheets" section
# (the Output paragraph ending with "and taken to their attendance sheet.")
$rng = $d.Content
$found = $rng.Find.Execute("and taken to their attendance sheet.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'and taken to their attendance sheet.'"
}

# Expand to the full paragraph (including its end-of-paragraph mark) so the
# new content is inserted as sibling paragraphs immediately after it.
$anchorPara = $rng.Paragraphs(1)
$insertAt = $d.Range($anchorPara.Range.Start, $anchorPara.Range.End)

$newParagraphsXml = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
    </w:rPr>
    <w:t>Edit Attendance Sheets</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Bold"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Description: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Admins can edit their attendance sheets</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Bold"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Actor: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>User</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>(Admins)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Bold"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Input: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">The user </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>adds new details editing the old details provided during creation of the attendance sheet.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:b/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Bold"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Output: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times-Roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Users will be displayed their new attendance sheet.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="1440"/>
    <w:rPr>
      <w:b/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
    </w:rPr>
  </w:pPr>
</w:p>

'@

$flatOpc = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
$newParagraphsXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$null = $insertAt.InsertXML($flatOpc, "After")
